# Recompute Student A's and Student C's scores after a Task 3 point bump,
# and keep the derived "Points"/"Score" summary sheets in sync.

$wb = $excel.ActiveWorkbook

$students = $wb.Worksheets.Item("Students")
$points   = $wb.Worksheets.Item("Points")

# --- Students sheet -------------------------------------------------------
# Student A (row 2): Task 3 score 2 -> 3, so Total 22 -> 23, grade 2- -> 2
$students.Range("G2").Value = 3
$students.Range("D2").Value = "23.0 (74.2%)"
# "2" looks numeric, so force text storage like the rest of the Grade column.
$students.Range("B2").NumberFormat = "@"
$students.Range("B2").Value = "2"
$students.Range("B2").Style = "Normal"

# Student C (row 4): Task 3 score 4 -> 5, so Total 17 -> 18, grade stays 3-
$students.Range("G4").Value = 5
$students.Range("D4").Value = "18.0 (58.1%)"

# Normalize the "Total" column formatting on the rest of the roster to the
# "<points>.0 (<pct>%)" form used above.
$students.Range("D3").Value  = "27.0 (87.1%)"
$students.Range("D5").Value  = "9.0 (29.0%)"
$students.Range("D6").Value  = "6.0 (19.4%)"
$students.Range("D7").Value  = "27.0 (87.1%)"
$students.Range("D8").Value  = "13.0 (41.9%)"
$students.Range("D9").Value  = "19.0 (61.3%)"
$students.Range("D10").Value = "20.0 (64.5%)"
$students.Range("D11").Value = "22.0 (71.0%)"

# --- Points sheet -----------------------------------------------------------
# Student C moves from 17 points to 18 points (row 5 stays the same student
# group/position, just the point total changes).
$points.Range("A5").Value = 18

# Student A splits off the "22 points" group (which now only contains
# Student J) into its own new "23 points" row, inserted before the existing
# 27-point row.
$points.Range("C8").Value = 1
$points.Range("D8").Value = "Student J"

$points.Rows("9").Insert()

# New row 9 points value (A) keeps the bold/bordered/centered look used for
# the rest of the "Points" column.
$points.Range("A9").Value = 23
$points.Range("A9").Font.Bold = $true
$points.Range("A9").Borders.LineStyle = 1
$points.Range("A9").HorizontalAlignment = -4108
$points.Range("A9").VerticalAlignment = -4160

# Grade "2" looks numeric, so force text storage the way the rest of the
# Grade column (e.g. "3", "6") is stored.
$points.Range("B9").NumberFormat = "@"
$points.Range("B9").Value = "2"
$points.Range("B9").Style = "Normal"

$points.Range("C9").Value = 1
$points.Range("D9").Value = "Student A"
